$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 256
$ws.Range("B5").Value = 0.1
$ws.Range("B6").Value = 50
$ws.Range("B7").Value = 3125

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
